$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new value for D4 (res.company.due_cost_service_id row) - new requirements string
$ws.Range("D4").Value = "l10n_it_ricevute_bancarie"

# Update the selected/active cell to A19 as in the saved view state
$ws.Range("A19").Select()
